$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "xxxxxxx" placeholder matriculation numbers in A3:A5 with the
# real matriculation numbers.
$ws.Range("A3").Value = 40399626
$ws.Range("A4").Value = 40410790
$ws.Range("A5").Value = 40399574


# Update the active selection to A5, matching the saved cursor position.
$ws.Range("A5").Select()
